$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("output1Salted")

# Update the B column (y-values) with the new salted/smoothed output data
$newValues = New-Object 'object[,]' 100,1
$newValues[0,0] = -14
$newValues[1,0] = -5
$newValues[2,0] = -9
$newValues[3,0] = 6
$newValues[4,0] = 3
$newValues[5,0] = 69
$newValues[6,0] = 13
$newValues[7,0] = 113
$newValues[8,0] = 122
$newValues[9,0] = 137
$newValues[10,0] = 87
$newValues[11,0] = 105
$newValues[12,0] = 216
$newValues[13,0] = 165
$newValues[14,0] = 183
$newValues[15,0] = 275
$newValues[16,0] = 275
$newValues[17,0] = 311
$newValues[18,0] = 323
$newValues[19,0] = 434
$newValues[20,0] = 483
$newValues[21,0] = 515
$newValues[22,0] = 518
$newValues[23,0] = 549
$newValues[24,0] = 664
$newValues[25,0] = 694
$newValues[26,0] = 691
$newValues[27,0] = 763
$newValues[28,0] = 800
$newValues[29,0] = 857
$newValues[30,0] = 994
$newValues[31,0] = 1063
$newValues[32,0] = 1130
$newValues[33,0] = 1182
$newValues[34,0] = 1251
$newValues[35,0] = 1260
$newValues[36,0] = 1389
$newValues[37,0] = 1417
$newValues[38,0] = 1473
$newValues[39,0] = 1621
$newValues[40,0] = 1659
$newValues[41,0] = 1726
$newValues[42,0] = 1881
$newValues[43,0] = 1961
$newValues[44,0] = 2059
$newValues[45,0] = 2071
$newValues[46,0] = 2165
$newValues[47,0] = 2352
$newValues[48,0] = 2377
$newValues[49,0] = 2519
$newValues[50,0] = 2642
$newValues[51,0] = 2694
$newValues[52,0] = 2850
$newValues[53,0] = 2951
$newValues[54,0] = 2991
$newValues[55,0] = 3089
$newValues[56,0] = 3298
$newValues[57,0] = 3340
$newValues[58,0] = 3496
$newValues[59,0] = 3620
$newValues[60,0] = 3760
$newValues[61,0] = 3877
$newValues[62,0] = 3987
$newValues[63,0] = 4112
$newValues[64,0] = 4256
$newValues[65,0] = 4384
$newValues[66,0] = 4463
$newValues[67,0] = 4668
$newValues[68,0] = 4734
$newValues[69,0] = 4858
$newValues[70,0] = 5070
$newValues[71,0] = 5147
$newValues[72,0] = 5294
$newValues[73,0] = 5515
$newValues[74,0] = 5587
$newValues[75,0] = 5742
$newValues[76,0] = 5947
$newValues[77,0] = 6115
$newValues[78,0] = 6193
$newValues[79,0] = 6450
$newValues[80,0] = 6586
$newValues[81,0] = 6758
$newValues[82,0] = 6932
$newValues[83,0] = 7083
$newValues[84,0] = 7273
$newValues[85,0] = 7369
$newValues[86,0] = 7605
$newValues[87,0] = 7724
$newValues[88,0] = 7901
$newValues[89,0] = 8061
$newValues[90,0] = 8327
$newValues[91,0] = 8513
$newValues[92,0] = 8691
$newValues[93,0] = 8885
$newValues[94,0] = 9039
$newValues[95,0] = 9203
$newValues[96,0] = 9429
$newValues[97,0] = 9656
$newValues[98,0] = 9781
$newValues[99,0] = 10024
$ws.Range("B1:B100").Value = $newValues

# Reposition the chart on the sheet (it was nudged down/right slightly)
$chartObj = $ws.ChartObjects().Item(1)
$chartObj.Left = 146.275
$chartObj.Top = 16.2
$chartObj.Width = 443.5
$chartObj.Height = 216.0

# Reposition the chart title (user dragged it to a custom spot)
$chart = $chartObj.Chart
$title = $chart.ChartTitle
$title.Left = 137.08146544181977
$title.Top = 8.0
